$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44 (shifts existing rows 44..162 down to 45..163)
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with the latest weekly record
$ws.Cells.Item(44,1).Value  = 5
$ws.Cells.Item(44,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(44,3).Value  = "Maule"
$ws.Cells.Item(44,4).Value  = 45251
$ws.Cells.Item(44,5).Value  = 7
$ws.Cells.Item(44,6).Value  = 100112022
$ws.Cells.Item(44,7).Value  = "Arveja Verde"
$ws.Cells.Item(44,8).Value  = "Sin especificar"
$ws.Cells.Item(44,9).Value  = "Primera"
$ws.Cells.Item(44,10).Value = 500
$ws.Cells.Item(44,11).Value = 25000
$ws.Cells.Item(44,12).Value = 25000
$ws.Cells.Item(44,13).Value = 25000
$ws.Cells.Item(44,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(44,15).Value = "Región del Maule"
$ws.Cells.Item(44,16).Value = 1000
$ws.Cells.Item(44,17).Value = 25
$ws.Cells.Item(44,18).Value = "Hortaliza"
